$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("hpi")
$ws1.Range("B2").Value = "These skin manifestations are highly specific for CREST syndrome, indicating the presence of the disease."
$ws1.Range("D2").Value = "If the patient does not report other gastrointestinal symptoms, it may suggest that Type 2 Achalasia is not present."
$ws1.Range("B3").Value = "Dysphagia is a common symptom associated with Type 2 Achalasia, which can occur in patients with CREST syndrome."
$ws1.Range("B4").Value = "Raynaud's phenomenon is a classic feature of CREST syndrome, supporting the diagnosis."
$ws1.Range("C4").Value = "Normal esophageal motility on manometry"
$ws1.Range("D4").Value = "Normal findings on esophageal motility testing would argue against the presence of Type 2 Achalasia."
$ws1.Range("B5").Value = "Reflux symptoms can be indicative of esophageal motility disorders like Type 2 Achalasia, which is associated with CREST syndrome."
$ws1.Range("C5").Value = "Absence of digital ulcers or calcinosis"
$ws1.Range("D5").Value = "These findings are often associated with CREST syndrome; their absence weakens the likelihood of the diagnosis."
$ws1.Range("B6").Value = "Pulmonary hypertension is a known complication of CREST syndrome, supporting the diagnosis."
$ws1.Range("D6").Value = "A lack of family history may suggest a lower risk for developing CREST syndrome."

$ws2 = $wb.Worksheets.Item("hist")
$ws2.Range("D2").Value = "The absence of skin changes or calcinosis makes CREST syndrome less likely, as these are common features."
$ws2.Range("C3").Value = "No history of digital ulcers"
$ws2.Range("D3").Value = "Digital ulcers are a common complication in CREST syndrome; their absence suggests a lower likelihood of the diagnosis."
$ws2.Range("C4").Value = "No history of telangiectasia"
$ws2.Range("D4").Value = "Telangiectasia is a hallmark of CREST syndrome; its absence weakens the case for this diagnosis."
$ws2.Range("C5").Value = "No previous diagnosis of systemic sclerosis"
$ws2.Range("D5").Value = "Systemic sclerosis is a broader condition that encompasses CREST syndrome; its absence suggests a lower likelihood of CREST."
$ws2.Range("B6").Value = "GERD can be a complication of esophageal motility disorders like achalasia, which is relevant in the context of CREST syndrome."
$ws2.Range("C6").Value = "No history of dysphagia unrelated to achalasia"
$ws2.Range("D6").Value = "Dysphagia that is not related to achalasia may indicate other causes, reducing the likelihood of CREST syndrome."

$ws3 = $wb.Worksheets.Item("soc")
$ws3.Range("B2").Value = "A family history of autoimmune diseases can indicate a genetic predisposition to conditions like CREST syndrome, which is an autoimmune disorder."
$ws3.Range("D2").Value = "Absence of a family history of CREST syndrome or related autoimmune diseases suggests a lower likelihood of the diagnosis."
$ws3.Range("B3").Value = "Certain occupations may increase exposure to environmental toxins, which have been associated with autoimmune diseases, including CREST syndrome."
$ws3.Range("C3").Value = "Engagement in healthy lifestyle practices"
$ws3.Range("D3").Value = "A healthy lifestyle, including regular exercise and a balanced diet, may reduce the risk of developing autoimmune conditions, including CREST syndrome."
$ws3.Range("B4").Value = "Raynaud's phenomenon is a common symptom in CREST syndrome, and a social history indicating this condition supports the diagnosis."
$ws3.Range("C4").Value = "No history of smoking or substance abuse"
$ws3.Range("D4").Value = "Lack of smoking or substance abuse is associated with a lower risk of developing autoimmune diseases, which may argue against the diagnosis."
$ws3.Range("B5").Value = "Chronic stress can exacerbate autoimmune conditions, and participation in high-stress activities may correlate with the development of CREST syndrome."
$ws3.Range("C5").Value = "Stable mental health with no history of anxiety or depression"
$ws3.Range("D5").Value = "Stable mental health may indicate lower stress levels, which can be a contributing factor to autoimmune conditions, thus arguing against CREST syndrome."
$ws3.Range("A6").Value = "History of gastrointestinal issues in family members"
$ws3.Range("B6").Value = "A family history of gastrointestinal issues may suggest a genetic link to conditions like Type 2 Achalasia, which can occur alongside CREST syndrome."
$ws3.Range("C6").Value = "No significant gastrointestinal complaints"
$ws3.Range("D6").Value = "Absence of gastrointestinal complaints suggests a lower likelihood of Type 2 Achalasia, which is often associated with CREST syndrome."

$ws4 = $wb.Worksheets.Item("obj")
$ws4.Range("C2").Value = "Normal vital signs"
$ws4.Range("D2").Value = "Stable vital signs may indicate the absence of significant systemic involvement, which is less common in CREST syndrome."
$ws4.Range("B3").Value = "This is a common symptom associated with CREST syndrome, indicating vascular involvement."
$ws4.Range("C3").Value = "Absence of skin changes"
$ws4.Range("D3").Value = "Lack of skin thickening or other sclerodermatous changes suggests that CREST syndrome is unlikely."
$ws4.Range("B4").Value = "This symptom is indicative of esophageal motility disorders, which are common in Type 2 Achalasia."
$ws4.Range("C4").Value = "Normal respiratory exam"
$ws4.Range("D4").Value = "A normal respiratory exam would argue against pulmonary complications often seen in CREST syndrome."
$ws4.Range("A5").Value = "Pulmonary hypertension"
$ws4.Range("B5").Value = "This can occur in patients with CREST syndrome and is a significant complication that can be detected on physical exam."
$ws4.Range("C5").Value = "No signs of esophageal dilation"
$ws4.Range("D5").Value = "The absence of esophageal dilation on physical exam would suggest that Type 2 Achalasia is not present."
$ws4.Range("A6").Value = "Telangiectasia"
$ws4.Range("B6").Value = "The presence of these small dilated blood vessels is a common finding in CREST syndrome."
$ws4.Range("C6").Value = "Normal capillary refill"
$ws4.Range("D6").Value = "Normal capillary refill time suggests adequate peripheral circulation, which may argue against Raynaud's phenomenon."

$ws5 = $wb.Worksheets.Item("test")
$ws5.Range("C2").Value = "Normal esophageal motility studies"
$ws5.Range("D2").Value = "Normal motility studies would suggest that Type 2 Achalasia is not present, arguing against the diagnosis."
$ws5.Range("D3").Value = "A negative result for these antibodies, which are associated with systemic sclerosis, would argue against CREST syndrome."
$ws5.Range("B4").Value = "This imaging finding is indicative of Type 2 Achalasia, which can occur in conjunction with CREST syndrome."
$ws5.Range("C4").Value = "Normal chest X-ray"
$ws5.Range("D4").Value = "A normal chest X-ray would not show any signs of esophageal dilation or other abnormalities associated with Type 2 Achalasia."
$ws5.Range("A5").Value = "Presence of calcinosis on imaging studies"
$ws5.Range("B5").Value = "Calcinosis is a common manifestation of CREST syndrome, supporting the diagnosis."
$ws5.Range("C5").Value = "Absence of skin changes or telangiectasia"
$ws5.Range("D5").Value = "The absence of these features would suggest that CREST syndrome is not present."
$ws5.Range("A6").Value = "Elevated anti-centromere antibodies"
$ws5.Range("B6").Value = "These antibodies are specifically associated with CREST syndrome, providing strong support for the diagnosis."
$ws5.Range("C6").Value = "Negative results for other autoimmune markers"
$ws5.Range("D6").Value = "Negative results for markers commonly associated with autoimmune diseases would argue against the presence of CREST syndrome."
